$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price column to remain text (it already stores text like "46.000.55")
# by switching to a text number format before writing, then restoring the
# original (default) style so no visible formatting change is introduced.
$priceRange = $ws.Range("D2:D51")
$origStyle = $priceRange.Style
$priceRange.NumberFormat = "@"

$ws.Range("D2").Value = "46.166.46"
$ws.Range("D3").Value = "2.454.04"
$ws.Range("D4").Value = "0.998"
$ws.Range("D5").Value = "322.59"
$ws.Range("D6").Value = "105.45"
$ws.Range("D9").Value = "0.540"
$ws.Range("D10").Value = "36.22"
$ws.Range("D11").Value = "0.0808"
$ws.Range("D12").Value = "0.123"
$ws.Range("D13").Value = "18.48"
$ws.Range("D14").Value = "7.09"
$ws.Range("D15").Value = "2.833.11"
$ws.Range("D16").Value = "2.452.58"
$ws.Range("D17").Value = "0.844"
$ws.Range("D18").Value = "45.941.81"
$ws.Range("D19").Value = "12.58"
$ws.Range("D20").Value = "6.45"
$ws.Range("D21").Value = "0.0₃0934"
$ws.Range("D22").Value = "71.77"
$ws.Range("D23").Value = "2.38"
$ws.Range("D24").Value = "248.17"
$ws.Range("D25").Value = "2.52"
$ws.Range("D26").Value = "26.07"
$ws.Range("D29").Value = "9.69"
$ws.Range("D30").Value = "33.70"
$ws.Range("D31").Value = "49.51"
$ws.Range("D32").Value = "0.130"
$ws.Range("D33").Value = "20.48"
$ws.Range("D36").Value = "0.0764"
$ws.Range("D37").Value = "4.57"
$ws.Range("D39").Value = "2.93"
$ws.Range("D40").Value = "127.19"
$ws.Range("D43").Value = "21.08"
$ws.Range("D44").Value = "0.0293"
$ws.Range("D45").Value = "1.964.48"
$ws.Range("D46").Value = "2.98"
$ws.Range("D47").Value = "2.12"
$ws.Range("D49").Value = "9.20"
$ws.Range("D50").Value = "77.88"
$ws.Range("D51").Value = "4.90"

$priceRange.Style = $origStyle

# Volume(1h) column is already text (leading/trailing spaces + % sign keep it text).
$ws.Range("E2").Value = "  +3.60%  "
$ws.Range("E3").Value = "  +0.86%  "
$ws.Range("E4").Value = "  -0.10%  "
$ws.Range("E5").Value = "  +3.31%  "
$ws.Range("E6").Value = "  +3.56%  "
$ws.Range("E7").Value = "  +0.98%  "
$ws.Range("E8").Value = "  -0.13%  "
$ws.Range("E9").Value = "  +6.10%  "
$ws.Range("E10").Value = "  +2.33%  "
$ws.Range("E11").Value = "  +0.95%  "
$ws.Range("E12").Value = "  -0.95%  "
$ws.Range("E13").Value = "  -1.65%  "
$ws.Range("E14").Value = "  +1.80%  "
$ws.Range("E15").Value = "  +0.69%  "
$ws.Range("E16").Value = "  +0.90%  "
$ws.Range("E17").Value = "  +0.86%  "
$ws.Range("E18").Value = "  +3.32%  "
$ws.Range("E19").Value = "  +0.66%  "
$ws.Range("E20").Value = "  +0.44%  "
$ws.Range("E21").Value = "  +2.63%  "
$ws.Range("E22").Value = "  +4.14%  "
$ws.Range("E23").Value = "  +2.12%  "
$ws.Range("E24").Value = "  +2.83%  "
$ws.Range("E25").Value = "  +0.96%  "
$ws.Range("E26").Value = "  +3.36%  "
$ws.Range("E28").Value = "  -4.71%  "
$ws.Range("E29").Value = "  -0.15%  "
$ws.Range("E30").Value = "  +1.50%  "
$ws.Range("E31").Value = "  +1.85%  "
$ws.Range("E32").Value = "  +6.59%  "
$ws.Range("E33").Value = "  +4.96%  "
$ws.Range("E34").Value = "  +1.83%  "
$ws.Range("E35").Value = "  -0.19%  "
$ws.Range("E36").Value = "  -0.04%  "
$ws.Range("E37").Value = "  +1.62%  "
$ws.Range("E38").Value = "  -0.37%  "
$ws.Range("E39").Value = "  +0.98%  "
$ws.Range("E40").Value = "  +0.87%  "
$ws.Range("E41").Value = "  +0.51%  "
$ws.Range("E42").Value = "  +1.69%  "
$ws.Range("E43").Value = "  -3.80%  "
$ws.Range("E44").Value = "  +1.02%  "
$ws.Range("E45").Value = "  +0.77%  "
$ws.Range("E46").Value = "  +1.36%  "
$ws.Range("E47").Value = "  -2.70%  "
$ws.Range("E48").Value = "  +10.30%  "
$ws.Range("E49").Value = "  -4.09%  "
$ws.Range("E50").Value = "  +5.23%  "
$ws.Range("E51").Value = "  +5.95%  "
